$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.365.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.805.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.02%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'315.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.97%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.15%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5490"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.33%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3854"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.75%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.44%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'42.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.22%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +3.66%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.188"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.15%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.347"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.59%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.802.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'92.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.41%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.06438"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'17.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.83%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.68%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'28.370.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.125"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'158.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.64%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.70%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.396"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.009.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.24%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'123.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.26%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.126"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.1018"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.42%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.743"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.22%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.672"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.88%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.2330"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +15.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.06364"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +4.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'8.842"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +10.81%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'11.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9991"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.44%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.39%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'13.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.65%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5977"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.20%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.686"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.88%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'124.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.16%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.79%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.148"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.61%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06905"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.17%  "
$ws.Range("E51").Style = "Normal"
